# Updated .dta files, improved literate programming cells
#
# var_details.xlsx keeps one row per Stata variable (position, name, type,
# isnumeric, format, vallab, varlab) for each questionnaire/instrument.
# The underlying .dta files were refreshed, which renamed a few variables
# and their formats/labels on the "director_data" sheet, and added a new
# "duplicatecheck" variable row to the "cso_data" sheet (mirroring the one
# that already exists on "director_data").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# director_data: variable renames / formats on existing rows
# ---------------------------------------------------------------------
$director = $wb.Worksheets.Item("director_data")

# row 3: zone -> subcounty
$director.Range("B3").Value = "subcounty"
$director.Range("C3").Value = "str38"
$director.Range("E3").Value = "%38s"
$director.Range("G3").Value = "Subcounty"

# row 8: format widened
$director.Range("C8").Value = "str562"
$director.Range("E8").Value = "%562s"

# row 10: format widened
$director.Range("C10").Value = "str127"
$director.Range("E10").Value = "%127s"

# row 47: duplicatecheck format changed from str1707 to strL
$director.Range("C47").Value = "strL"
$director.Range("E47").Value = "%9s"

# ---------------------------------------------------------------------
# cso_data: same variable rename as director_data row 3, plus a new
# duplicatecheck row inserted right before the trailing instrument row
# ---------------------------------------------------------------------
$cso = $wb.Worksheets.Item("cso_data")

$cso.Range("C3").Value = "str19"
$cso.Range("E3").Value = "%19s"

# insert a new row 60 (pushes the old row 60 "instrument" row down to 61)
$cso.Range("A60:G60").Insert()

$cso.Range("A60").Value = 59
$cso.Range("B60").Value = "duplicatecheck"
$cso.Range("C60").Value = "strL"
$cso.Range("D60").Value = 0
$cso.Range("E60").Value = "%9s"
$cso.Range("F60").Value = ""
$cso.Range("G60").Value = "duplicate check"

# the "position" counter for the row that got pushed down (old row 60,
# the trailing "instrument" row) advances by one, from 59 to 60
$cso.Range("A61").Value = 60
